$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the hidden "Skill Heatmap" sheet to "Calculations" and unhide it
# ---------------------------------------------------------------------------
$calc = $wb.Worksheets.Item("Skill Heatmap")
$calc.Name = "Calculations"
$calc.Visible = -1

# ---------------------------------------------------------------------------
# 2. Re-order sheets: move "Calculations" to the end (after
#    "Planning and Stabilizing Teams"), so tab order becomes:
#    Readme, Settings, Survey Sheet, Input and results,
#    Planning and Stabilizing Teams, Calculations
# ---------------------------------------------------------------------------
$planning = $wb.Worksheets.Item("Planning and Stabilizing Teams")
$calc.Move($null, $planning)

# ---------------------------------------------------------------------------
# 3. Update the "Calculations" totals column (R2:R18) so a row with no
#    skill entered (sum = 0) shows blank instead of 0.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 18; $r++) {
    $calc.Range("R$r").Formula = "=IF(SUM(F$r`:Q$r)>0,SUM(F$r`:Q$r),`"`")"
}

# ---------------------------------------------------------------------------
# 4. Selection on the Calculations sheet moves to R23
# ---------------------------------------------------------------------------
$calc.Range("R23").Select()

# ---------------------------------------------------------------------------
# 5. Restore selections on other sheets (unchanged in the diff, but make
#    sure nothing drifted)
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")
$settings.Range("D28").Select()

# ---------------------------------------------------------------------------
# 6. Readme becomes the active/selected tab; window scroll resets
# ---------------------------------------------------------------------------
$readme = $wb.Worksheets.Item("Readme")
$readme.Activate()
$readme.Range("F35").Select()

$excel.ActiveWindow.WindowState = -4143  # xlNormal (no-op safeguard)
$wb.Windows.Item(1).ScrollColumn = 1
$wb.Windows.Item(1).ScrollRow = 1
